$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 877.55554
$ws.Range("J29").Value = 1100
$ws.Range("L29").Value = 3300
$ws.Range("N29").Value = -3862
$ws.Range("H62").Value = 27295
$ws.Range("I62").Value = 4279.2
$ws.Range("K62").Value = 4279.2
$ws.Range("M62").Value = -3655.2
$ws.Range("H65").Value = 27295
$ws.Range("I65").Value = 4279.2
$ws.Range("K65").Value = 21396
$ws.Range("M65").Value = -18276
$ws.Range("H88").Value = 1707.4445
$ws.Range("J88").Value = 1909.5714
$ws.Range("L88").Value = 1909.5714
$ws.Range("N88").Value = -2721.5714
$ws.Range("H91").Value = 1707.4445
$ws.Range("J91").Value = 1909.5714
$ws.Range("L91").Value = 1909.5714
$ws.Range("N91").Value = -4717.5714
$ws.Range("H97").Value = 1998.3334
$ws.Range("J97").Value = 1998.3334
$ws.Range("L97").Value = 5995.0002
$ws.Range("N97").Value = -6987.0002
$ws.Range("H98").Value = 3651.195
$ws.Range("I98").Value = 3340.457
$ws.Range("K98").Value = 3340.457
$ws.Range("M98").Value = -1842.457
$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 10000
$ws.Range("N105").Value = -16988
$ws.Range("H122").Value = 3651.195
$ws.Range("I122").Value = 3340.457
$ws.Range("K122").Value = 10021.371
$ws.Range("M122").Value = -7571.370999999999
$ws.Range("H135").Value = 398.78262
$ws.Range("I135").Value = 392.2381
$ws.Range("K135").Value = 3530.1429
$ws.Range("M135").Value = -995.1428999999998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 273.5
$ws.Range("I25").Value = 273.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 273.5
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 128.5
$ws.Range("N25").ClearContents()
$ws.Range("H32").Value = 8020.4653
$ws.Range("I32").Value = 5248.7754
$ws.Range("J32").Value = 23110.777
$ws.Range("K32").Value = 5248.7754
$ws.Range("L32").Value = 23110.777
$ws.Range("M32").Value = -4961.7754
$ws.Range("N32").Value = -23684.777
$ws.Range("H97").Value = 957.05
$ws.Range("I97").Value = 924.71875
$ws.Range("K97").Value = 924.71875
$ws.Range("M97").Value = -428.71875
$ws.Range("H110").Value = 25001566
$ws.Range("I110").Value = 29413084
$ws.Range("K110").Value = 29413084
$ws.Range("M110").Value = -29411039
$ws.Range("H122").Value = 2724.3484
$ws.Range("I122").Value = 2314.4424
$ws.Range("K122").Value = 6943.3272
$ws.Range("M122").Value = -4493.3272
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 54995
$ws.Range("J124").Value = 54995
$ws.Range("L124").Value = 54995
$ws.Range("N124").Value = -64815
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2802.2632
$ws.Range("I31").Value = 2059.182
$ws.Range("K31").Value = 2059.182
$ws.Range("M31").Value = -1764.182
$ws.Range("H34").Value = 2802.2632
$ws.Range("I34").Value = 2059.182
$ws.Range("K34").Value = 2059.182
$ws.Range("M34").Value = -1857.182
$ws.Range("H62").Value = 52796.6
$ws.Range("I62").Value = 2495.5
$ws.Range("J62").Value = 86330.664
$ws.Range("K62").Value = 2495.5
$ws.Range("L62").Value = 86330.664
$ws.Range("M62").Value = -1871.5
$ws.Range("N62").Value = -87578.664
$ws.Range("H65").Value = 52796.6
$ws.Range("I65").Value = 2495.5
$ws.Range("J65").Value = 86330.664
$ws.Range("K65").Value = 12477.5
$ws.Range("L65").Value = 431653.32
$ws.Range("M65").Value = -9357.5
$ws.Range("N65").Value = -437893.32
$ws.Range("H107").Value = 545.73914
$ws.Range("I107").Value = 531.9
$ws.Range("K107").Value = 531.9
$ws.Range("M107").Value = 1388.1
$ws.Range("H134").Value = 3847.1333
$ws.Range("I134").Value = 3170.25
$ws.Range("J134").Value = 6554.6665
$ws.Range("K134").Value = 9510.75
$ws.Range("L134").Value = 19663.9995
$ws.Range("M134").Value = -6975.75
$ws.Range("N134").Value = -24733.9995
$ws.Range("H141").Value = 427599.6
$ws.Range("J141").Value = 427599.6
$ws.Range("L141").Value = 427599.6
$ws.Range("N141").Value = -437959.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4011.6667
$ws.Range("J76").Value = 4011.6667
$ws.Range("L76").Value = 12035.0001
$ws.Range("N76").Value = -12801.0001
$ws.Range("H79").Value = 4011.6667
$ws.Range("J79").Value = 4011.6667
$ws.Range("L79").Value = 12035.0001
$ws.Range("N79").Value = -14687.0001
$ws.Range("H98").Value = 724.25
$ws.Range("I98").Value = 632.3333
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 1896.9999
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -398.9999
$ws.Range("N98").Value = -5996
$ws.Range("H107").Value = 625.63635
$ws.Range("J107").Value = 377
$ws.Range("L107").Value = 1131
$ws.Range("N107").Value = -4971
$ws.Range("H132").Value = 2722.9539
$ws.Range("I132").Value = 1166.6666
$ws.Range("K132").Value = 10499.9994
$ws.Range("M132").Value = -7969.999400000001
$ws.Range("H134").Value = 2741.7856
$ws.Range("I134").Value = 1806.5385
$ws.Range("K134").Value = 5419.6155
$ws.Range("M134").Value = -349.6154999999999
$ws.Range("H137").Value = 3038.88
$ws.Range("I137").Value = 2636.2
$ws.Range("J137").Value = 3139.55
$ws.Range("K137").Value = 7908.599999999999
$ws.Range("L137").Value = 9418.650000000001
$ws.Range("M137").Value = -2808.599999999999
$ws.Range("N137").Value = -19618.65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 31666.666
$ws.Range("J52").Value = 31666.666
$ws.Range("L52").Value = 31666.666
$ws.Range("N52").Value = -32184.666
$ws.Range("H97").Value = 751.34485
$ws.Range("I97").Value = 780.2941
$ws.Range("J97").Value = 710.3333
$ws.Range("K97").Value = 780.2941
$ws.Range("L97").Value = 710.3333
$ws.Range("M97").Value = -284.2941
$ws.Range("N97").Value = -1702.3333
$ws.Range("H104").Value = 27500
$ws.Range("J104").Value = 27500
$ws.Range("L104").Value = 27500
$ws.Range("N104").Value = -34488
$ws.Range("H134").Value = 71775.336
$ws.Range("J134").Value = 71775.336
$ws.Range("L134").Value = 215326.008
$ws.Range("N134").Value = -220396.008
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 725.2
$ws.Range("I58").Value = 725.2
$ws.Range("K58").Value = 725.2
$ws.Range("M58").Value = -465.2
$ws.Range("H61").Value = 7099.222
$ws.Range("J61").Value = 2495
$ws.Range("L61").Value = 2495
$ws.Range("N61").Value = -2899
$ws.Range("H82").Value = 2446.125
$ws.Range("I82").Value = 2192.5
$ws.Range("J82").Value = 2699.75
$ws.Range("K82").Value = 2192.5
$ws.Range("L82").Value = 2699.75
$ws.Range("M82").Value = -1831.5
$ws.Range("N82").Value = -3421.75
$ws.Range("H85").Value = 2446.125
$ws.Range("I85").Value = 2192.5
$ws.Range("J85").Value = 2699.75
$ws.Range("K85").Value = 2192.5
$ws.Range("L85").Value = 2699.75
$ws.Range("M85").Value = -944.5
$ws.Range("N85").Value = -5195.75
$ws.Range("H100").Value = 3207.3076
$ws.Range("I100").Value = 2911.875
$ws.Range("K100").Value = 2911.875
$ws.Range("M100").Value = -2370.875
$ws.Range("H113").Value = 7099.222
$ws.Range("J113").Value = 2495
$ws.Range("L113").Value = 2495
$ws.Range("N113").Value = -6835
$ws.Range("H127").Value = 49950
$ws.Range("J127").Value = 49950
$ws.Range("L127").Value = 49950
$ws.Range("N127").Value = -59870
$ws.Range("H135").Value = 45500
$ws.Range("J135").Value = 45500
$ws.Range("L135").Value = 45500
$ws.Range("N135").Value = -55640
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2858327.8
$ws.Range("I3").Value = 6667133.5
$ws.Range("K3").Value = 6667133.5
$ws.Range("M3").Value = -6667019.5
$ws.Range("H81").Value = 9528447
$ws.Range("I81").Value = 1849.9
$ws.Range("K81").Value = 3699.8
$ws.Range("M81").Value = -2638.8
$ws.Range("H84").Value = 9528447
$ws.Range("I84").Value = 1849.9
$ws.Range("K84").Value = 18499
$ws.Range("M84").Value = -13195
$ws.Range("H93").Value = 56856.145
$ws.Range("I93").Value = 38999.5
$ws.Range("J93").Value = 63998.8
$ws.Range("K93").Value = 38999.5
$ws.Range("L93").Value = 63998.8
$ws.Range("M93").Value = -36503.5
$ws.Range("N93").Value = -68990.8
$ws.Range("H122").Value = 2366
$ws.Range("I122").Value = 2318.182
$ws.Range("K122").Value = 6954.545999999999
$ws.Range("M122").Value = -4504.545999999999
$ws.Range("H132").Value = 369044.1
$ws.Range("I132").Value = 576563.9399999999
$ws.Range("J132").Value = 5884.4
$ws.Range("K132").Value = 1729691.82
$ws.Range("L132").Value = 17653.2
$ws.Range("M132").Value = -1727161.82
$ws.Range("N132").Value = -22713.2
